# Actualización automática 2025-10-17 14:30:08
#
# This script applies the monthly sales-update edit to the workbook:
#  - "VENTAS POR GRUPO"     : per-client / per-product-group raw sales
#  - "VENTA MENSUAL"        : per-client monthly sales roll-up (octubre = col F)
#  - "CUMPLIMIENTO MENSUAL" : per-group budget-vs-actual roll-up
#
# Two clients received new "octubre" sales that used to be 0 (or, for one of
# them, the "octubre" figure simply grew):
#   - row 30 "INTRIAGO ALVARADO BRENDA ALEJANDRA" -> PORCELANATO group: +1590.22
#   - row 47 "SALAZAR BALLADARES MARIA ANGELICA"  -> 240X80 PORCELANATO group: +356.16
#
# All of the other touched cells are downstream roll-ups of those two raw
# values, so we write every one of them explicitly (the workbook stores only
# literal values, no formulas).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": raw per-client / per-group sales figures
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# INTRIAGO ALVARADO BRENDA ALEJANDRA - PORCELANATO group
$wsGrupo.Range("M30").Value = 1590.22

# SALAZAR BALLADARES MARIA ANGELICA - 240X80 PORCELANATO group
$wsGrupo.Range("D47").Value = 356.16

# Row 59 totals: count of clients with a positive value in each group
# ("X de 57"). Both PORCELANATO (M) and 240X80 PORCELANATO (D) gained one
# more client with a positive amount.
$wsGrupo.Range("D59").Value = "4 de 57"
$wsGrupo.Range("M59").Value = "6 de 57"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": per-client monthly roll-up (column F = octubre)
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F30").Value = 1590.22
$wsMensual.Range("F47").Value = 1558.83
$wsMensual.Range("F59").Value = 18359.84

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": budget vs. actual per product group
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3: 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 2794.18
$wsCumpl.Range("E3").Value = 17593.2974217135
$wsCumpl.Range("F3").Value = 0.1370537385377597

# Row 12: PORCELANATO (combined)
$wsCumpl.Range("D12").Value = 7854.93
$wsCumpl.Range("E12").Value = 40769.13
$wsCumpl.Range("F12").Value = 0.1615440997728285

# Row 14: TOTAL
$wsCumpl.Range("D14").Value = 21928.46
$wsCumpl.Range("E14").Value = 77969.53284188786
$wsCumpl.Range("F14").Value = 0.2195085143973509
